# Apply scheduled profit-recalculation updates to the Seraph_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) produced by the runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 851
$ws.Range("I2").Value = 586.25
$ws.Range("J2").Value = 1027.5
$ws.Range("K2").Value = 586.25
$ws.Range("L2").Value = 1027.5
$ws.Range("M2").Value = -473.25
$ws.Range("N2").Value = -1253.5
$ws.Range("H34").Value = 2828.8333
$ws.Range("I34").Value = 2828.8333
$ws.Range("K34").Value = 2828.8333
$ws.Range("M34").Value = -2625.8333
$ws.Range("H36").Value = 2828.8333
$ws.Range("I36").Value = 2828.8333
$ws.Range("K36").Value = 2828.8333
$ws.Range("M36").Value = -2113.8333
$ws.Range("H43").Value = 17499
$ws.Range("I43").Value = 17499
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 17499
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -17430
$ws.Range("H107").Value = 1032
$ws.Range("I107").Value = 1291.3684
$ws.Range("J107").Value = 328
$ws.Range("K107").Value = 1291.3684
$ws.Range("L107").Value = 328
$ws.Range("M107").Value = 628.6315999999999
$ws.Range("N107").Value = -4168
$ws.Range("H132").Value = 2487
$ws.Range("I132").Value = 1528.4166
$ws.Range("K132").Value = 4585.2498
$ws.Range("M132").Value = -2055.2498
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1638.7368
$ws.Range("I45").Value = 1337.4286
$ws.Range("K45").Value = 1337.4286
$ws.Range("M45").Value = -960.4286
$ws.Range("H63").Value = 5319
$ws.Range("I63").Value = 1108.8334
$ws.Range("J63").Value = 8927.714
$ws.Range("K63").Value = 1108.8334
$ws.Range("L63").Value = 8927.714
$ws.Range("M63").Value = -422.8334
$ws.Range("N63").Value = -10299.714
$ws.Range("H66").Value = 5319
$ws.Range("I66").Value = 1108.8334
$ws.Range("J66").Value = 8927.714
$ws.Range("K66").Value = 5544.166999999999
$ws.Range("L66").Value = 44638.57
$ws.Range("M66").Value = -2112.166999999999
$ws.Range("N66").Value = -51502.57
$ws.Range("H102").Value = 1950
$ws.Range("I102").Value = 2183.3333
$ws.Range("J102").Value = 1250
$ws.Range("K102").Value = 2183.3333
$ws.Range("L102").Value = 1250
$ws.Range("M102").Value = -561.3332999999998
$ws.Range("N102").Value = -4494
$ws.Range("H132").Value = 2952.6667
$ws.Range("I132").Value = 2946.75
$ws.Range("K132").Value = 8840.25
$ws.Range("M132").Value = -6310.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5559.2144
$ws.Range("I105").Value = 4448.4614
$ws.Range("K105").Value = 4448.4614
$ws.Range("M105").Value = -2701.4614
$ws.Range("H134").Value = 2829.6667
$ws.Range("I134").Value = 2835.6
$ws.Range("K134").Value = 8506.799999999999
$ws.Range("M134").Value = -5971.799999999999
$ws.Range("H135").Value = 60998
$ws.Range("J135").Value = 60998
$ws.Range("L135").Value = 60998
$ws.Range("N135").Value = -71138

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 93555.69500000001
$ws.Range("I22").Value = 133277.67
$ws.Range("J22").Value = 4181.25
$ws.Range("K22").Value = 133277.67
$ws.Range("L22").Value = 4181.25
$ws.Range("M22").Value = -132927.67
$ws.Range("N22").Value = -4881.25
$ws.Range("H25").Value = 1764
$ws.Range("I25").Value = 1092.5714
$ws.Range("K25").Value = 1092.5714
$ws.Range("M25").Value = -918.5714
$ws.Range("H86").Value = 8526.5
$ws.Range("I86").Value = 8468.666999999999
$ws.Range("J86").Value = 8700
$ws.Range("K86").Value = 8468.666999999999
$ws.Range("L86").Value = 8700
$ws.Range("M86").Value = -7345.666999999999
$ws.Range("N86").Value = -10946
$ws.Range("H89").Value = 8526.5
$ws.Range("I89").Value = 8468.666999999999
$ws.Range("J89").Value = 8700
$ws.Range("K89").Value = 42343.335
$ws.Range("L89").Value = 43500
$ws.Range("M89").Value = -36727.335
$ws.Range("N89").Value = -54732
$ws.Range("H122").Value = 2845.25
$ws.Range("I122").Value = 1555.4
$ws.Range("J122").Value = 4995
$ws.Range("K122").Value = 4666.200000000001
$ws.Range("L122").Value = 14985
$ws.Range("M122").Value = -2216.200000000001
$ws.Range("N122").Value = -19885

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 738.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 738.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2215.5
$ws.Range("N68").Value = -3837.5
$ws.Range("H71").Value = 738.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 738.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 6646.5
$ws.Range("N71").Value = -14758.5
$ws.Range("H103").Value = 224.16667
$ws.Range("J103").Value = 293
$ws.Range("L103").Value = 879
$ws.Range("N103").Value = -2637
$ws.Range("H131").Value = 4462.8423
$ws.Range("J131").Value = 5905.615
$ws.Range("L131").Value = 17716.845
$ws.Range("N131").Value = -27796.845
$ws.Range("H132").Value = 2873.375
$ws.Range("I132").Value = 2799.8
$ws.Range("J132").Value = 2996
$ws.Range("K132").Value = 25198.2
$ws.Range("L132").Value = 26964
$ws.Range("M132").Value = -22668.2
$ws.Range("N132").Value = -32024
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1069.4667
$ws.Range("I97").Value = 734
$ws.Range("J97").Value = 3250
$ws.Range("K97").Value = 734
$ws.Range("L97").Value = 3250
$ws.Range("M97").Value = -238
$ws.Range("N97").Value = -4242
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("H132").Value = 3387
$ws.Range("I132").Value = 3330.5
$ws.Range("K132").Value = 9991.5
$ws.Range("M132").Value = -7461.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 15000
$ws.Range("I16").Value = 15000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -14830
$ws.Range("H22").Value = 1293.2222
$ws.Range("I22").Value = 1079.6666
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 1079.6666
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -784.6666
$ws.Range("N22").Value = -1990
$ws.Range("H27").Value = 1293.2222
$ws.Range("I27").Value = 1079.6666
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 1079.6666
$ws.Range("L27").Value = 1400
$ws.Range("M27").Value = -972.6666
$ws.Range("N27").Value = -1614
$ws.Range("H61").Value = 3081.875
$ws.Range("I61").Value = 1912.5
$ws.Range("J61").Value = 4251.25
$ws.Range("K61").Value = 1912.5
$ws.Range("L61").Value = 4251.25
$ws.Range("M61").Value = -1710.5
$ws.Range("N61").Value = -4655.25
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("H113").Value = 3081.875
$ws.Range("I113").Value = 1912.5
$ws.Range("J113").Value = 4251.25
$ws.Range("K113").Value = 1912.5
$ws.Range("L113").Value = 4251.25
$ws.Range("M113").Value = 257.5
$ws.Range("N113").Value = -8591.25
$ws.Range("H122").Value = 3461.0322
$ws.Range("I122").Value = 2888.6667
$ws.Range("J122").Value = 3695.182
$ws.Range("K122").Value = 8666.000100000001
$ws.Range("L122").Value = 11085.546
$ws.Range("M122").Value = -6216.000100000001
$ws.Range("N122").Value = -15985.546
$ws.Range("H136").Value = 5575
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 5150
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 15450
$ws.Range("M136").Value = -15450
$ws.Range("N136").Value = -20550
$ws.Range("N16").ClearContents()
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 50000
$ws.Range("I48").Value = 50000
$ws.Range("K48").Value = 50000
$ws.Range("M48").Value = -49431
$ws.Range("H81").Value = 1429.8
$ws.Range("I81").Value = 1429.8
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2859.6
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1798.6
$ws.Range("H84").Value = 1429.8
$ws.Range("I84").Value = 1429.8
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 14298
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -8994
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("H122").Value = 1602.4286
$ws.Range("I122").Value = 1470.3529
$ws.Range("J122").Value = 2163.75
$ws.Range("K122").Value = 4411.0587
$ws.Range("L122").Value = 6491.25
$ws.Range("M122").Value = -1961.0587
$ws.Range("N122").Value = -11391.25
$ws.Range("H136").Value = 4964.4287
$ws.Range("I136").Value = 4964.4287
$ws.Range("K136").Value = 14893.2861
$ws.Range("M136").Value = -12343.2861
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("N116").ClearContents()
